$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last data row (old row 29, CIRO FERNANDO CASELLES BENAVIDES / period 2507)
# carries the special "closing" border formatting for the table. We want that
# same closing formatting to end up on the new last row once CIRO's row and
# JHON's row are removed, so copy its format onto the row that will become
# the new last row (old row 28 - DONIS CENTENO MARTINEZ / period 2502)
# before any rows shift.
$ws.Range("B29:J29").Copy() | Out-Null
$ws.Range("B28:J28").PasteSpecial(-4122, -4142, $false, $false) | Out-Null

# Remove JHON JAIRO HERNANDEZ MEJIA's single row entirely.
$ws.Rows("16").Delete()

# Remove CIRO FERNANDO CASELLES BENAVIDES's single row entirely (shifted up
# by one after the previous delete, from 29 to 28).
$ws.Rows("28").Delete()

# Refresh the summary figures now that only two workers remain (LUIS CARLOS
# CENTENO DIAZ and DONIS CENTENO MARTINEZ), six periods each at 52000.
$ws.Range("E11").Value = 624000
$ws.Range("C13").Value = 2

# The "Nombre Trabajador" column no longer needs to fit the long
# "CIRO FERNANDO CASELLES BENAVIDES" string, so it shrinks back down.
$ws.Columns("D").AutoFit() | Out-Null
